# Updated RAD to add MD Central Registration Number to 2 tax forms.
# (This workbook records the Katalon test-run log: each row's "Result"
# (col A) / "Date" (col B) cells get stamped after a fresh test pass.)

$wb = $excel.ActiveWorkbook

# --- Estimated sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Estimated")
$ws.Cells.Item(2, 2).Value = "Fri Oct 06 12:43:01 EDT 2023"
$ws.Cells.Item(3, 2).Value = "Fri Oct 06 12:43:20 EDT 2023"

$ws.Cells.Item(4, 1).Value = "Pass"
$ws.Cells.Item(4, 2).Value = "Fri Oct 06 12:43:39 EDT 2023"
$ws.Cells.Item(5, 1).Value = "Pass"
$ws.Cells.Item(5, 2).Value = "Fri Oct 06 12:43:57 EDT 2023"
$ws.Cells.Item(6, 1).Value = "Pass"
$ws.Cells.Item(6, 2).Value = "Fri Oct 06 12:44:16 EDT 2023"
$ws.Cells.Item(7, 1).Value = "Pass"
$ws.Cells.Item(7, 2).Value = "Fri Oct 06 12:44:36 EDT 2023"

# --- Existing sheet -----------------------------------------------------
$ws = $wb.Worksheets.Item("Existing")
$ws.Cells.Item(2, 2).Value = "Fri Oct 06 12:44:55 EDT 2023"
$ws.Cells.Item(3, 2).Value = "Fri Oct 06 12:45:12 EDT 2023"
$ws.Cells.Item(4, 2).Value = "Fri Oct 06 12:45:30 EDT 2023"
$ws.Cells.Item(5, 2).Value = "Fri Oct 06 12:45:47 EDT 2023"
$ws.Cells.Item(6, 2).Value = "Fri Oct 06 12:46:05 EDT 2023"
$ws.Cells.Item(7, 2).Value = "Fri Oct 06 12:46:22 EDT 2023"
$ws.Cells.Item(8, 2).Value = "Fri Oct 06 12:46:40 EDT 2023"
$ws.Cells.Item(9, 2).Value = "Fri Oct 06 12:46:57 EDT 2023"
$ws.Cells.Item(10, 2).Value = "Fri Oct 06 12:47:14 EDT 2023"
$ws.Cells.Item(11, 2).Value = "Fri Oct 06 12:47:32 EDT 2023"
$ws.Cells.Item(12, 2).Value = "Fri Oct 06 12:47:49 EDT 2023"

# --- Extension sheet ------------------------------------------------------
$ws = $wb.Worksheets.Item("Extension")
$ws.Cells.Item(2, 2).Value = "Fri Oct 06 12:48:07 EDT 2023"
$ws.Cells.Item(3, 2).Value = "Fri Oct 06 12:48:26 EDT 2023"
$ws.Cells.Item(4, 2).Value = "Fri Oct 06 12:48:46 EDT 2023"
$ws.Cells.Item(5, 2).Value = "Fri Oct 06 12:49:04 EDT 2023"
$ws.Cells.Item(6, 2).Value = "Fri Oct 06 12:49:21 EDT 2023"
$ws.Cells.Item(7, 2).Value = "Fri Oct 06 12:49:39 EDT 2023"

# --- NewTaxReturn sheet ----------------------------------------------------
$ws = $wb.Worksheets.Item("NewTaxReturn")
$ws.Cells.Item(2, 2).Value = "Fri Oct 06 12:49:58 EDT 2023"
$ws.Cells.Item(3, 2).Value = "Fri Oct 06 12:50:17 EDT 2023"
$ws.Cells.Item(4, 2).Value = "Fri Oct 06 12:50:35 EDT 2023"
$ws.Cells.Item(5, 2).Value = "Fri Oct 06 12:50:53 EDT 2023"
$ws.Cells.Item(6, 2).Value = "Fri Oct 06 12:51:12 EDT 2023"
$ws.Cells.Item(7, 2).Value = "Fri Oct 06 12:51:30 EDT 2023"
$ws.Cells.Item(8, 2).Value = "Fri Oct 06 12:51:48 EDT 2023"
$ws.Cells.Item(9, 2).Value = "Fri Oct 06 12:52:07 EDT 2023"
$ws.Cells.Item(10, 2).Value = "Fri Oct 06 12:52:25 EDT 2023"
$ws.Cells.Item(11, 2).Value = "Fri Oct 06 12:52:44 EDT 2023"
$ws.Cells.Item(12, 2).Value = "Fri Oct 06 12:53:02 EDT 2023"
$ws.Cells.Item(13, 2).Value = "Fri Oct 06 12:53:20 EDT 2023"
$ws.Cells.Item(14, 2).Value = "Fri Oct 06 12:53:39 EDT 2023"
$ws.Cells.Item(15, 2).Value = "Fri Oct 06 12:53:57 EDT 2023"
$ws.Cells.Item(16, 2).Value = "Fri Oct 06 12:54:16 EDT 2023"

# --- Personal_EL sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_EL")
$ws.Cells.Item(2, 2).Value = "Fri Oct 06 12:54:34 EDT 2023"

# --- Personal_IND sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_IND")
$ws.Cells.Item(2, 2).Value = "Fri Oct 06 12:54:57 EDT 2023"
$ws.Cells.Item(3, 2).Value = "Fri Oct 06 12:55:16 EDT 2023"
$ws.Cells.Item(4, 2).Value = "Fri Oct 06 12:55:35 EDT 2023"
$ws.Cells.Item(5, 2).Value = "Fri Oct 06 12:55:54 EDT 2023"
$ws.Cells.Item(6, 2).Value = "Fri Oct 06 12:56:13 EDT 2023"

# --- Personal_JNT sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Personal_JNT")
$ws.Cells.Item(2, 2).Value = "Fri Oct 06 12:56:33 EDT 2023"
$ws.Cells.Item(3, 2).Value = "Fri Oct 06 12:57:00 EDT 2023"
$ws.Cells.Item(4, 2).Value = "Fri Oct 06 12:57:26 EDT 2023"
$ws.Cells.Item(5, 2).Value = "Fri Oct 06 12:57:54 EDT 2023"
$ws.Cells.Item(6, 2).Value = "Fri Oct 06 12:58:21 EDT 2023"
